$d = $word.ActiveDocument
Write-Host "paragraphs count:" $d.Paragraphs.Count
